$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 34
$ws1.Range("F3").Value = 182
$ws1.Range("F6").Value = 559
$ws1.Range("F7").Value = 1744
$ws1.Range("F11").Value = 1902
$ws1.Range("F14").Value = 439
$ws1.Range("F15").Value = 8
$ws1.Range("F16").Value = 276
$ws1.Range("F18").Value = 6
$ws1.Range("F23").Value = 1027
$ws1.Range("F26").Value = 173
$ws1.Range("F27").Value = 255
$ws1.Range("F28").Value = 279

# Sheet "全部类型" (All types) - update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 34
$ws4.Range("F3").Value = 182
$ws4.Range("F6").Value = 559
$ws4.Range("F7").Value = 1744
$ws4.Range("F12").Value = 1902
$ws4.Range("F15").Value = 439
$ws4.Range("F16").Value = 8
$ws4.Range("F17").Value = 276
$ws4.Range("F19").Value = 6
$ws4.Range("F24").Value = 1027
$ws4.Range("F27").Value = 173
$ws4.Range("F28").Value = 255
$ws4.Range("F29").Value = 279
